$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")
$ws.Range("B9").Value = "test"
